$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("provincia") is being recurated from a dimension (sdmx-dimension:refArea /
# dim / URI-Provincia) to a measure (iaest-measure:provincia / medida / xsd:int).
$ws.Range("G2").Value = "iaest-measure:provincia"
$ws.Range("G3").Value = "medida"
$ws.Range("G4").Value = "xsd:int"
